# Rename sheets to match the updated JSON-converted vocabulary naming.
$wb = $excel.ActiveWorkbook

$wsTechnique = $wb.Worksheets.Item("Measurement technique")
$wsTechnique.Name = "Technique"

$wsProperty = $wb.Worksheets.Item("Measured property")
$wsProperty.Name = "Measured property #parameter"

# Keep "Technique" as the active tab (matches original authoring state),
# while moving the cursor/selection on the "Measured property #parameter"
# sheet to E60.
[void]$wsTechnique.Activate()
[void]$wsProperty.Range("E60").Select()
[void]$wsTechnique.Activate()
